# Update countries & provincias Spain
# Refresh the case/death statistics for the affected country rows on the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7819134
$ws.Range("C4").Value = 42271
$ws.Range("D4").Value = 5011191
$ws.Range("E4").Value = 2590455
$ws.Range("G4").Value = 707
$ws.Range("H4").Value = 217488

# Row 5 - India
$ws.Range("B5").Value = 6903806
$ws.Range("C5").Value = 70818
$ws.Range("D5").Value = 5903170
$ws.Range("E5").Value = 894082
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 106554

# Row 26 - Alemania
$ws.Range("B26").Value = 315454
$ws.Range("C26").Value = 4341
$ws.Range("E26").Value = 38087
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 9667

# Row 27 - Israel
$ws.Range("B27").Value = 285336
$ws.Range("C27").Value = 3855
$ws.Range("D27").Value = 221571
$ws.Range("E27").Value = 61901

# Row 89 - Costa de Marfil
$ws.Range("B89").Value = 19982
$ws.Range("C89").Value = 47
$ws.Range("D89").Value = 19626
$ws.Range("E89").Value = 236

# Row 101 - Namibia
$ws.Range("B101").Value = 11781
$ws.Range("C101").Value = 67
$ws.Range("D101").Value = 9759
$ws.Range("E101").Value = 1895
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 127

# Row 114 - Zimbabue
$ws.Range("B114").Value = 7951
$ws.Range("C114").Value = 32
$ws.Range("D114").Value = 6446
$ws.Range("E114").Value = 1276

# Row 118 - Cabo Verde
$ws.Range("B118").Value = 6717
$ws.Range("C118").Value = 93
$ws.Range("D118").Value = 5821
$ws.Range("E118").Value = 825

# Row 121 - Malaui
$ws.Range("B121").Value = 5809
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 4626
$ws.Range("E121").Value = 1003

# Row 132 - Ruanda
$ws.Range("B132").Value = 4885
$ws.Range("C132").Value = 2
$ws.Range("D132").Value = 3542
$ws.Range("E132").Value = 1314

# Row 135 - Siria
$ws.Range("B135").Value = 4566
$ws.Range("C135").Value = 62
$ws.Range("D135").Value = 1212
$ws.Range("E135").Value = 3139
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 215

# Row 147 - Mali
$ws.Range("B147").Value = 3235
$ws.Range("C147").Value = 25
$ws.Range("D147").Value = 2506
$ws.Range("E147").Value = 598

# Row 148 - Botsuana
$ws.Range("B148").Value = 3219
$ws.Range("C148").Value = 47
$ws.Range("E148").Value = 2367

# Row 161 - Togo
$ws.Range("B161").Value = 1907
$ws.Range("C161").Value = 9
$ws.Range("D161").Value = 1426
$ws.Range("E161").Value = 432

# Row 179 - Comoras
$ws.Range("B179").Value = 495
$ws.Range("C179").Value = 4
$ws.Range("D179").Value = 475
$ws.Range("E179").Value = 13

# Row 192 - Bermudas
$ws.Range("B192").Value = 182
$ws.Range("C192").Value = 1
$ws.Range("E192").Value = 3
